# Actualización automática de tasas-transfi.xlsx
# Updates the "conversion del dia" text block on Hoja1 and the
# rate values on the "tasas" sheet.

$wb = $excel.ActiveWorkbook

# --- Hoja1: update the daily conversion summary text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.61 = 5847.65 pesos`n✅ 5847.65 pesos = 1.6 = 940.97 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- tasas: update N10/O10 and N12/O12 values ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 622.481
$wsTasas.Range("O10").Value = 3640.05

$wsTasas.Range("N12").Value = 3659
$wsTasas.Range("O12").Value = 588.787
